$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# maj de la BOM : la ligne 6 (bobine/condensateur) est maintenant disponible au garage
$ws.Range("E6").Value = "oui"

# Met à jour la cellule active / la sélection affichée dans le classeur
$ws.Range("C7").Select()
